$wb = $excel.ActiveWorkbook

$wsCounts = $wb.Worksheets.Item("st counts & ridges")
$wsYarn   = $wb.Worksheets.Item("yarn proportions")

# ---------------------------------------------------------------------------
# "yarn proportions" sheet: core numeric edit + formula tweaks
# ---------------------------------------------------------------------------

# Cast-on count for MC (column C) drops from 2 to 1 skein/strand ratio.
$wsYarn.Range("C2").Value = 1

# Simplify the short-row ratio formulas (drop the redundant "*B_" multiply).
$wsYarn.Range("D12").Formula = "=B3/B2"
$wsYarn.Range("D13").Formula = "=B4/B3"
$wsYarn.Range("D14").Formula = "=B5/B4"
$wsYarn.Range("D15").Formula = "=B6/B5"

# ---------------------------------------------------------------------------
# Selections / active sheet: the workbook was left with "yarn proportions"
# as the active tab, a new selection there, and a new selection back on
# "st counts & ridges".
# ---------------------------------------------------------------------------

$wsCounts.Activate()
$wsCounts.Range("G43").Select()

$wsYarn.Activate()
$wsYarn.Range("E14").Select()

$wb.RecalculateFullRebuild()
